$wb = $excel.ActiveWorkbook

# ---------------------------------------------------------------------
# Sheet "App-draft": update the existing request row and clone it into a
# new second request row (row 3), then adjust the sheet view/selection.
# ---------------------------------------------------------------------
$wsApp = $wb.Worksheets.Item("App-draft")
$wsApp.Activate()

$wsApp.Range("M1").Value = "Allocate"

$wsApp.Range("D2").Value = "Choose Sunday, December 1st, 2024"
$wsApp.Range("K2").Value = "Choose Sunday, December 1st, 2024"
$wsApp.Range("L2").Value = "Choose Sunday, December 1st, 2024"

$wsApp.Range("A2:M2").Copy($wsApp.Range("A3:M3"))
$wsApp.Range("E3").Value = "Kapila"
$wsApp.Range("F3").Value = "786567865v"
$wsApp.Range("D3").Value = "Choose Sunday, December 1st, 2024"
$wsApp.Range("K3").Value = "Choose Sunday, December 1st, 2024"
$wsApp.Range("L3").Value = "Choose Sunday, December 1st, 2024"

$wsApp.Range("L3").Select()

# ---------------------------------------------------------------------
# Sheet "Section3": the complaint date moves forward.
# ---------------------------------------------------------------------
$wsSection3 = $wb.Worksheets.Item("Section3")
$wsSection3.Activate()
$wsSection3.Range("K2").Value = "Choose Sunday, December 1st, 2024"
$wsSection3.Range("K2").Select()

# ---------------------------------------------------------------------
# Sheet "chequeDetails": new voucher/cheque numbers and issued date.
# ---------------------------------------------------------------------
$wsCheque = $wb.Worksheets.Item("chequeDetails")
$wsCheque.Activate()
$wsCheque.Range("A2").Value = 8977787
$wsCheque.Range("B2").Value = 877676
$wsCheque.Range("C2").Value = "Choose Sunday, December 1st, 2024"
$wsCheque.Range("B2").Select()

# ---------------------------------------------------------------------
# New sheet "chequePrint" appended after "chequeDetails".
# ---------------------------------------------------------------------
$lastSheet = $wb.Worksheets.Item($wb.Worksheets.Count)
$wsPrint = $wb.Worksheets.Add($null, $lastSheet)
$wsPrint.Name = "chequePrint"

$wsPrint.Columns.Item(1).ColumnWidth = 11.5
$wsPrint.Columns.Item(2).ColumnWidth = 28.83
$wsPrint.Columns.Item(3).ColumnWidth = 15.17

$wsPrint.Range("A1").Value = "Language"
$wsPrint.Range("B1").Value = "Payment Date"
$wsPrint.Range("C1").Value = "Reference No"
$wsPrint.Range("D1").Value = "Officer"
$wsPrint.Range("A1:D1").Interior.ThemeColor = 3
$wsPrint.Range("A1:D1").Interior.TintAndShade = 0.39997558519241921
$wsPrint.Range("A1:D1").Font.Bold = $true

$wsPrint.Range("A2").Value = "Sinhala"
$wsPrint.Range("B2").Value = "Choose Sunday, December 1st, 2024"
$wsPrint.Range("B2").NumberFormat = "m/d/yyyy"
$wsPrint.Range("C2").Value = 678995
$wsPrint.Range("D2").Value = "Cheque Officer"

$wsPrint.Range("B7").Select()

# ---------------------------------------------------------------------
# Restore "chequeDetails" as the active sheet/tab.
# ---------------------------------------------------------------------
$wsCheque.Activate()
